# Updated main GSC export files
# - Chart: append 4 new daily rows (2025-11-19 .. 2025-11-22)
# - Critical issues: refresh reason/source/validation/pages breakdown

$wb = $excel.ActiveWorkbook

# Helper: write a literal text value into a cell without letting the
# date-autodetection turn "2025-11-19"-like strings into date serials.
function Set-TextCell($sheet, $row, $col, $text) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# Chart sheet: four new rows of daily totals
# ---------------------------------------------------------------------
$chart = $wb.Worksheets.Item("Chart")

Set-TextCell $chart 47 1 "2025-11-19"
$chart.Cells.Item(47, 2).Value = 127
$chart.Cells.Item(47, 3).Value = 213
$chart.Cells.Item(47, 4).Value = 24

Set-TextCell $chart 48 1 "2025-11-20"
$chart.Cells.Item(48, 2).Value = 127
$chart.Cells.Item(48, 3).Value = 213
$chart.Cells.Item(48, 4).Value = 26

Set-TextCell $chart 49 1 "2025-11-21"
$chart.Cells.Item(49, 2).Value = 127
$chart.Cells.Item(49, 3).Value = 213
$chart.Cells.Item(49, 4).Value = 38

Set-TextCell $chart 50 1 "2025-11-22"
$chart.Cells.Item(50, 2).Value = 127
$chart.Cells.Item(50, 3).Value = 213
$chart.Cells.Item(50, 4).Value = 24

# ---------------------------------------------------------------------
# Critical issues sheet: refreshed coverage breakdown (header unchanged)
# ---------------------------------------------------------------------
$critical = $wb.Worksheets.Item("Critical issues")

$criticalRows = @(
    @("Alternate page with proper canonical tag", "Website", "Failed", 40),
    @("Not found (404)", "Website", "Failed", 21),
    @("Duplicate, Google chose different canonical than user", "Google systems", "Failed", 36),
    @("Excluded by ‘noindex’ tag", "Website", "Not Started", 15),
    @("Server error (5xx)", "Website", "Not Started", 1),
    @("Blocked by robots.txt", "Website", "Not Started", 1),
    @("Page with redirect", "Website", "Started", 2),
    @("Crawled - currently not indexed", "Google systems", "Started", 9),
    @("Discovered - currently not indexed", "Google systems", "Started", 2)
)

$r = 2
foreach ($row in $criticalRows) {
    $critical.Cells.Item($r, 1).Value = $row[0]
    $critical.Cells.Item($r, 2).Value = $row[1]
    $critical.Cells.Item($r, 3).Value = $row[2]
    $critical.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
